$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.27784
$ws.Range("H2").Value = 30.83352
$ws.Range("I2").Value = 0.230301226653591
$ws.Range("J2").Value = 0.230301226653591
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 47.57896333333333
$ws.Range("N2").Value = 142.73689
$ws.Range("O2").Value = 0.450188452948237
$ws.Range("P2").Value = 0.4501884529482371
$ws.Range("Q2").Value = 489.0089725058666
$ws.Range("R2").Value = 4401.0807525528
$ws.Range("S2").Value = 0.1036789529392614
$ws.Range("T2").Value = 0.1036789529392614

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.27784
$ws.Range("H3").Value = 30.83352
$ws.Range("I3").Value = 0.230301226653591
$ws.Range("J3").Value = 0.230301226653591
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.829723666666666
$ws.Range("N3").Value = 29.489171
$ws.Range("O3").Value = 0.09300808131111737
$ws.Range("P3").Value = 0.09300808131111739
$ws.Range("Q3").Value = 101.0283270902133
$ws.Range("R3").Value = 909.2549438119199
$ws.Range("S3").Value = 0.02141987521464726
$ws.Range("T3").Value = 0.02141987521464727

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.27784
$ws.Range("H4").Value = 30.83352
$ws.Range("I4").Value = 0.230301226653591
$ws.Range("J4").Value = 0.230301226653591
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.933664
$ws.Range("N4").Value = 32.800992
$ws.Range("O4").Value = 0.1034534789405002
$ws.Range("P4").Value = 0.1034534789405003
$ws.Range("Q4").Value = 112.37444920576
$ws.Range("R4").Value = 1011.37004285184
$ws.Range("S4").Value = 0.02382546310157865
$ws.Range("T4").Value = 0.02382546310157865

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 10.27784
$ws.Range("H5").Value = 30.83352
$ws.Range("I5").Value = 0.230301226653591
$ws.Range("J5").Value = 0.230301226653591
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 37.34441866666666
$ws.Range("N5").Value = 112.033256
$ws.Range("O5").Value = 0.3533499868001453
$ws.Range("P5").Value = 0.3533499868001453
$ws.Range("Q5").Value = 383.8199599490133
$ws.Range("R5").Value = 3454.37963954112
$ws.Range("S5").Value = 0.08137693539810366
$ws.Range("T5").Value = 0.08137693539810366

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 14.56812733333333
$ws.Range("H6").Value = 43.704382
$ws.Range("I6").Value = 0.3264360600001921
$ws.Range("J6").Value = 0.326436060000192
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 47.57896333333333
$ws.Range("N6").Value = 142.73689
$ws.Range("O6").Value = 0.450188452948237
$ws.Range("P6").Value = 0.4501884529482371
$ws.Range("Q6").Value = 693.1363962279976
$ws.Range("R6").Value = 6238.227566051979
$ws.Range("S6").Value = 0.1469577448380044
$ws.Range("T6").Value = 0.1469577448380043

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 14.56812733333333
$ws.Range("H7").Value = 43.704382
$ws.Range("I7").Value = 0.3264360600001921
$ws.Range("J7").Value = 0.326436060000192
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.829723666666666
$ws.Range("N7").Value = 29.489171
$ws.Range("O7").Value = 0.09300808131111737
$ws.Range("P7").Value = 0.09300808131111739
$ws.Range("Q7").Value = 143.2006660274802
$ws.Range("R7").Value = 1288.805994247322
$ws.Range("S7").Value = 0.03036119161137866
$ws.Range("T7").Value = 0.03036119161137866

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.56812733333333
$ws.Range("H8").Value = 43.704382
$ws.Range("I8").Value = 0.3264360600001921
$ws.Range("J8").Value = 0.326436060000192
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 10.933664
$ws.Range("N8").Value = 32.800992
$ws.Range("O8").Value = 0.1034534789405002
$ws.Range("P8").Value = 0.1034534789405003
$ws.Range("Q8").Value = 159.2830093718827
$ws.Range("R8").Value = 1433.547084346944
$ws.Range("S8").Value = 0.03377094605864975
$ws.Range("T8").Value = 0.03377094605864975

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.56812733333333
$ws.Range("H9").Value = 43.704382
$ws.Range("I9").Value = 0.3264360600001921
$ws.Range("J9").Value = 0.326436060000192
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 37.34441866666666
$ws.Range("N9").Value = 112.033256
$ws.Range("O9").Value = 0.3533499868001453
$ws.Range("P9").Value = 0.3533499868001453
$ws.Range("Q9").Value = 544.0382463253101
$ws.Range("R9").Value = 4896.344216927791
$ws.Range("S9").Value = 0.1153461774921593
$ws.Range("T9").Value = 0.1153461774921593

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.650137
$ws.Range("H10").Value = 7.950411
$ws.Range("I10").Value = 0.05938308067649115
$ws.Range("J10").Value = 0.05938308067649114
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 47.57896333333333
$ws.Range("N10").Value = 142.73689
$ws.Range("O10").Value = 0.450188452948237
$ws.Range("P10").Value = 0.4501884529482371
$ws.Range("Q10").Value = 126.09077115131
$ws.Range("R10").Value = 1134.81694036179
$ws.Range("S10").Value = 0.0267335772210499
$ws.Range("T10").Value = 0.0267335772210499

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.650137
$ws.Range("H11").Value = 7.950411
$ws.Range("I11").Value = 0.05938308067649115
$ws.Range("J11").Value = 0.05938308067649114
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 9.829723666666666
$ws.Range("N11").Value = 29.489171
$ws.Range("O11").Value = 0.09300808131111737
$ws.Range("P11").Value = 0.09300808131111739
$ws.Range("Q11").Value = 26.050114388809
$ws.Range("R11").Value = 234.451029499281
$ws.Range("S11").Value = 0.005523106396063731
$ws.Range("T11").Value = 0.005523106396063732

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.650137
$ws.Range("H12").Value = 7.950411
$ws.Range("I12").Value = 0.05938308067649115
$ws.Range("J12").Value = 0.05938308067649114
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.933664
$ws.Range("N12").Value = 32.800992
$ws.Range("O12").Value = 0.1034534789405002
$ws.Range("P12").Value = 0.1034534789405003
$ws.Range("Q12").Value = 28.975707511968
$ws.Range("R12").Value = 260.781367607712
$ws.Range("S12").Value = 0.006143386286187404
$ws.Range("T12").Value = 0.006143386286187404

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.650137
$ws.Range("H13").Value = 7.950411
$ws.Range("I13").Value = 0.05938308067649115
$ws.Range("J13").Value = 0.05938308067649114
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 37.34441866666666
$ws.Range("N13").Value = 112.033256
$ws.Range("O13").Value = 0.3533499868001453
$ws.Range("P13").Value = 0.3533499868001453
$ws.Range("Q13").Value = 98.96782565202399
$ws.Range("R13").Value = 890.7104308682159
$ws.Range("S13").Value = 0.02098301077319011
$ws.Range("T13").Value = 0.02098301077319011

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 17.13170833333333
$ws.Range("H14").Value = 51.395125
$ws.Range("I14").Value = 0.3838796326697257
$ws.Range("J14").Value = 0.3838796326697257
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 47.57896333333333
$ws.Range("N14").Value = 142.73689
$ws.Range("O14").Value = 0.450188452948237
$ws.Range("P14").Value = 0.4501884529482371
$ws.Range("Q14").Value = 815.1089226290276
$ws.Range("R14").Value = 7335.980303661249
$ws.Range("S14").Value = 0.1728181779499213
$ws.Range("T14").Value = 0.1728181779499213

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 17.13170833333333
$ws.Range("H15").Value = 51.395125
$ws.Range("I15").Value = 0.3838796326697257
$ws.Range("J15").Value = 0.3838796326697257
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 9.829723666666666
$ws.Range("N15").Value = 29.489171
$ws.Range("O15").Value = 0.09300808131111737
$ws.Range("P15").Value = 0.09300808131111739
$ws.Range("Q15").Value = 168.3999588545972
$ws.Range("R15").Value = 1515.599629691375
$ws.Range("S15").Value = 0.03570390808902772
$ws.Range("T15").Value = 0.03570390808902772

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 17.13170833333333
$ws.Range("H16").Value = 51.395125
$ws.Range("I16").Value = 0.3838796326697257
$ws.Range("J16").Value = 0.3838796326697257
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 10.933664
$ws.Range("N16").Value = 32.800992
$ws.Range("O16").Value = 0.1034534789405002
$ws.Range("P16").Value = 0.1034534789405003
$ws.Range("Q16").Value = 187.3123426626667
$ws.Range("R16").Value = 1685.811083964
$ws.Range("S16").Value = 0.03971368349408445
$ws.Range("T16").Value = 0.03971368349408445

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 17.13170833333333
$ws.Range("H17").Value = 51.395125
$ws.Range("I17").Value = 0.3838796326697257
$ws.Range("J17").Value = 0.3838796326697257
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 37.34441866666666
$ws.Range("N17").Value = 112.033256
$ws.Range("O17").Value = 0.3533499868001453
$ws.Range("P17").Value = 0.3533499868001453
$ws.Range("Q17").Value = 639.7736884752221
$ws.Range("R17").Value = 5757.963196277
$ws.Range("S17").Value = 0.1356438631366922
$ws.Range("T17").Value = 0.1356438631366922
